# Apply multiplication-problem text replacements to the document.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "898×8="; New = "428×7=" },
    @{ Old = "611×3="; New = "969×4=" },
    @{ Old = "338×8="; New = "810×9=" },
    @{ Old = "704×3="; New = "969×4=" },
    @{ Old = "823×5="; New = "851×5=" },
    @{ Old = "513×5="; New = "892×2=" },
    @{ Old = "652×7="; New = "852×4=" },
    @{ Old = "196×5="; New = "663×2=" },
    @{ Old = "712×8="; New = "196×9=" },
    @{ Old = "995×9="; New = "329×3=" },
    @{ Old = "881×5="; New = "976×7=" },
    @{ Old = "912×6="; New = "986×8=" },
    @{ Old = "175×6="; New = "138×9=" },
    @{ Old = "514×2="; New = "498×8=" },
    @{ Old = "184×2="; New = "810×5=" },
    @{ Old = "770×3="; New = "704×5=" },
    @{ Old = "817×7="; New = "155×2=" },
    @{ Old = "349×8="; New = "947×2=" },
    @{ Old = "180×2="; New = "393×3=" },
    @{ Old = "245×3="; New = "815×4=" },
    @{ Old = "904×4="; New = "642×8=" },
    @{ Old = "427×8="; New = "453×7=" },
    @{ Old = "569×9="; New = "291×2=" },
    @{ Old = "785×9="; New = "736×5=" },
    @{ Old = "454×3="; New = "636×6=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
